$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # validcreds
$ws2 = $wb.Worksheets.Item(2)   # invalidcreds
$ws3 = $wb.Worksheets.Item(3)   # Sheet3 -> "customer and project creds"

# --- Rename third sheet ---
$ws3.Name = "customer and project creds"

# --- invalidcreds (sheet2): update the failing-testcase credential data ---
# Order matters: it controls the order new shared strings are appended, so
# write these first, in the exact sequence that reproduces the target
# shared-string table ordering.
$ws2.Range("A3").Value = "ad"
$ws2.Range("B4").Value = "mana"
$ws2.Range("A5").Value = "min"
$ws2.Range("B5").Value = "admin"
$ws2.Range("A8").Value = "mana"
$ws2.Range("B8").Value = "ad"

# --- validcreds (sheet1): repurpose sheet for customer/project data ---
$ws1.Range("A1").Value = "Customer_Name"
$ws1.Range("B1").Value = "Project_Name"
$ws1.Range("A2").Value = "Yes Bank"
$ws1.Range("B2").Value = "Automate Web Application"

# --- customer and project creds (sheet3): same customer/project data ---
$ws3.Range("A1").Value = "Customer_Name"
$ws3.Range("B1").Value = "Project_Name"
$ws3.Range("A2").Value = "Yes Bank"
$ws3.Range("B2").Value = "Automate Web Application"

# Copy the header/value cell formatting from sheet1 onto sheet3 so the new
# sheet gets the same bold/yellow header style and bordered value style.
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:B1").PasteSpecial(-4122)
$ws1.Range("A2:B2").Copy()
$ws3.Range("A2:B2").PasteSpecial(-4122)

# --- Column widths ---
$ws1.Columns.Item(1).ColumnWidth = 17.5
$ws1.Columns.Item(2).ColumnWidth = 22.666666666666668
$ws3.Columns.Item(1).ColumnWidth = 15.833333333333334
$ws3.Columns.Item(2).ColumnWidth = 24.166666666666668

# --- Page setup for the new sheet ---
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# --- Selections ---
$ws1.Range("C8").Select()
$ws2.Range("B8").Select()
$ws3.Range("B10").Select()

# --- Active sheet/tab ---
$ws3.Activate()
